# "Update countries & provincias Spain" - refresh COVID-19 case counters.
#
# The sheet "Pais" is kept sorted by column B (Casos totales) descending.
# Refreshing a handful of countries' numbers changes their rank, which
# shuffles which country label sits on which row. Below we first write the
# updated numeric stats for the countries whose own row doesn't move
# (España, Brasil, Paises Bajos, Austria), then rewrite the small block of
# rows around Liberia/Sierra Leona whose relative order changes (Sierra
# Leona's count grows enough to overtake Liberia, Guadalupe, Birmania,
# Gibraltar and Brunei), and finally swap two same-valued ties
# (Namibia/San Vicente y las Granadinas and Comoras/San Pedro y Miquelon).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- España (row 5): totals refreshed, row position unchanged ---
$ws.Range("B5").Value = 245567
$ws.Range("C5").Value = 2579
$ws.Range("D5").Value = 146233
$ws.Range("E5").Value = 74234
$ws.Range("G5").Value = 276
$ws.Range("H5").Value = 25100

# --- Brasil (row 13): totals refreshed, row position unchanged ---
$ws.Range("B13").Value = 92630
$ws.Range("C13").Value = 521
$ws.Range("E13").Value = 48157
$ws.Range("G13").Value = 24
$ws.Range("H13").Value = 6434

# --- Paises Bajos (row 18): totals refreshed, row position unchanged ---
$ws.Range("B18").Value = 40236
$ws.Range("C18").Value = 445
$ws.Range("E18").Value = 34999
$ws.Range("G18").Value = 94
$ws.Range("H18").Value = 4987

# --- Austria (row 32): only active/recovered/critical columns move ---
$ws.Range("D32").Value = 13180
$ws.Range("E32").Value = 1782
$ws.Range("F32").Value = 114

# --- Rows 137-142: Sierra Leona's updated count (155) now outranks
#     Liberia/Guadalupe/Birmania/Gibraltar/Brunei, so each of those
#     countries' data shifts down one row and Sierra Leona's fresh data
#     takes row 137. Row 143 (Madagascar) is unaffected. ---
$ws.Range("A137").Value = "Sierra Leona"
$ws.Range("B137").Value = 155
$ws.Range("C137").Value = 19
$ws.Range("D137").Value = 21
$ws.Range("E137").Value = 126
$ws.Range("G137").Value = 1
$ws.Range("H137").Value = 8

$ws.Range("A138").Value = "Liberia"
$ws.Range("D138").Value = 45
$ws.Range("E138").Value = 89
$ws.Range("F138").Value = 0
$ws.Range("H138").Value = 18

$ws.Range("A139").Value = "Guadalupe"
$ws.Range("B139").Value = 152
$ws.Range("D139").Value = 95
$ws.Range("E139").Value = 45
$ws.Range("F139").Value = 6
$ws.Range("H139").Value = 12

$ws.Range("A140").Value = "Birmania"
$ws.Range("B140").Value = 151
$ws.Range("D140").Value = 31
$ws.Range("E140").Value = 114
$ws.Range("H140").Value = 6

$ws.Range("A141").Value = "Gibraltar"
$ws.Range("B141").Value = 144
$ws.Range("D141").Value = 131
$ws.Range("E141").Value = 13
$ws.Range("F141").Value = 0
$ws.Range("H141").Value = 0

$ws.Range("A142").Value = "Brunei"
$ws.Range("B142").Value = 138
$ws.Range("D142").Value = 126
$ws.Range("E142").Value = 11
$ws.Range("F142").Value = 2
$ws.Range("H142").Value = 1

# --- Rows 193/194: Namibia and San Vicente y las Granadinas are tied
#     (16/0/8/8/0/0/0), so only the labels swap. ---
$ws.Range("A193").Value = "San Vicente y las Granadinas"
$ws.Range("A194").Value = "Namibia"

# --- Rows 217/218: Comoras and San Pedro y Miquelon are tied
#     (1/0/0/1/0/0/0), so only the labels swap. ---
$ws.Range("A217").Value = "San Pedro y Miquelon"
$ws.Range("A218").Value = "Comoras"
